$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'60.405.10"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.93%  '
$ws.Range('D3').Value = "'2.612.62"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.96%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = "'585.19"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.25%  '
$ws.Range('D6').Value = "'143.31"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.31%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = "'0.597"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.04%  '
$ws.Range('D9').Value = "'6.52"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.17%  '
$ws.Range('E10').Value = '  -1.37%  '
$ws.Range('E11').Value = '  +1.00%  '
$ws.Range('D12').Value = "'0.373"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.98%  '
$ws.Range('D13').Value = "'3.077.12"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.77%  '
$ws.Range('D14').Value = "'24.71"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.62%  '
$ws.Range('D15').Value = "'60.404.68"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.90%  '
$ws.Range('E16').Value = '  -0.34%  '
$ws.Range('D17').Value = "'2.620.44"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.24%  '
$ws.Range('E18').Value = '  +0.79%  '
$ws.Range('D19').Value = "'4.64"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.07%  '
$ws.Range('D20').Value = "'346.60"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.36%  '
$ws.Range('D21').Value = "'6.91"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E22').Value = '  -0.73%  '
$ws.Range('D23').Value = "'0.535"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.38%  '
$ws.Range('D24').Value = "'63.52"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.18%  '
$ws.Range('E26').Value = '  -0.67%  '
$ws.Range('D27').Value = "'8.01"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.68%  '
$ws.Range('D28').Value = "'1.92"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.17%  '
$ws.Range('D29').Value = "'0.0₃0797"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.19%  '
$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D30').Value = "'6.44"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.70%  '
$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').Value = "'169.06"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.26%  '
$ws.Range('E32').Value = '  +0.00%  '
$ws.Range('D33').Value = "'19.48"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.20%  '
$ws.Range('D34').Value = "'1.01"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.19%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').Value = "'4.29"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = "'1.30"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +8.07%  '
$ws.Range('E37').Value = '  +3.38%  '
$ws.Range('D38').Value = "'319.05"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +6.57%  '
$ws.Range('E39').Value = '  +1.69%  '
$ws.Range('D40').Value = "'3.91"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.54%  '
$ws.Range('D41').Value = "'0.850"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.67%  '
$ws.Range('D42').Value = "'136.22"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.12%  '
$ws.Range('E43').Value = '  +0.30%  '
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').Value = "'0.610"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.36%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = "'19.92"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.37%  '
$ws.Range('D47').Value = "'5.03"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.62%  '
$ws.Range('D48').Value = "'0.0549"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.06%  '
$ws.Range('D49').Value = "'20.08"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.91%  '
$ws.Range('E50').Value = '  -0.27%  '
$ws.Range('D51').Value = "'10.74"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.43%  '
